# Online user tracking added
# Adds a new localization row (ESchoolLoginWarning) to the tr.xlsx resource sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 88

$ws.Range("A$newRow").Value = "ESchoolLoginWarning"
$ws.Range("B$newRow").Value = "E-Okul'dan veri cekilmesi esnasinda kullandiginiz kullanici bilgieri ile baska bir yerde oturum acmamaniz gerekmektedir. Oturum acildigi takdir de  uygulamadaki  veri aktarimi kesilecektir ve islem devam etmeyecektir."

# Carry over the formatting from the row above (instead of assigning a new
# Style object, which would mint a brand-new style/font combo) so the new
# row keeps the same cell styles already used throughout the sheet.
$ws.Range("A87").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122)

$ws.Range("B87").Copy()
$ws.Range("B$newRow").PasteSpecial(-4122)
